$wb = $excel.ActiveWorkbook

$wsRocket = $wb.Worksheets.Item("Rocket")
$wsManeuvers = $wb.Worksheets.Item("Maneuvers")

# --- Maneuvers sheet data edits ---
# Row 2: was "Plane Adjust" / Earth / 100 / 0 / 400 / 0
#        -> "Transfer to Lunar Orbit" / Earth / 3600 / 5310 / 0 / 0
$wsManeuvers.Range("A2").Value = "Transfer to Lunar Orbit"
$wsManeuvers.Range("B2").Value = "Earth"
$wsManeuvers.Range("C2").Value = 3600
$wsManeuvers.Range("D2").Value = 5310
$wsManeuvers.Range("E2").Value = 0
$wsManeuvers.Range("F2").Value = 0

# Row 3: was "Transfer A" / Earth / 101 / 4000 / 0 / 0
#        -> "Plane Adjust" / Earth / 7200 / 0 / -375 / 0
$wsManeuvers.Range("A3").Value = "Plane Adjust"
$wsManeuvers.Range("B3").Value = "Earth"
$wsManeuvers.Range("C3").Value = 7200
$wsManeuvers.Range("D3").Value = 0
$wsManeuvers.Range("E3").Value = -375
$wsManeuvers.Range("F3").Value = 0

# Row 4 (new): "Insertion" / Moon / 357000 / -850 / 0 / 0
$wsManeuvers.Range("A4").Value = "Insertion"
$wsManeuvers.Range("B4").Value = "Moon"
$wsManeuvers.Range("C4").Value = 357000
$wsManeuvers.Range("D4").Value = -850
$wsManeuvers.Range("E4").Value = 0
$wsManeuvers.Range("F4").Value = 0

# Row 5 (new): "Circularization" / Moon / 375000 / -270 / 0 / 0
$wsManeuvers.Range("A5").Value = "Circularization"
$wsManeuvers.Range("B5").Value = "Moon"
$wsManeuvers.Range("C5").Value = 375000
$wsManeuvers.Range("D5").Value = -270
$wsManeuvers.Range("E5").Value = 0
$wsManeuvers.Range("F5").Value = 0

# Column A width on Maneuvers sheet widened (closest achievable value)
$wsManeuvers.Columns.Item(1).ColumnWidth = 23

# --- Selections (cursor position) ---
$null = $wsRocket.Range("C24").Select()
$null = $wsManeuvers.Range("B20").Select()
